$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new rows before the existing row 520, shifting all following rows down.
$ws.Rows("520:523").Insert()

# Common values shared by the new rows (same product/market as surrounding rows).
$mercadoId = 9
$mercado = "Vega Central Mapocho de Santiago"
$region = "Metropolitana"
$fecha = 44706
$codreg = 13
$tipo = "Fruta"
$productoId = 100101
$producto = "Berries"
$categoriaId = 100101007
$categoria = "Kiwi"
$variedad = "Hayward"
$unidad = "$/bandeja 10 kilos"
$origen = "Provincia de Curicó"
$kgUnidad = 10

$newRows = @(
    @{ Row = 520; Calidad = "Especial";                 Volumen = 220; Precio = 10000; PrecioKg = 1000 },
    @{ Row = 521; Calidad = "Extra (doble especial)";    Volumen = 250; Precio = 12000; PrecioKg = 1200 },
    @{ Row = 522; Calidad = "Primera";                   Volumen = 280; Precio = 8000;  PrecioKg = 800 },
    @{ Row = 523; Calidad = "Segunda";                   Volumen = 250; Precio = 7000;  PrecioKg = 700 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = $mercadoId
    $ws.Cells.Item($row, 2).Value = $mercado
    $ws.Cells.Item($row, 3).Value = $region
    $ws.Cells.Item($row, 4).Value = $fecha
    $ws.Cells.Item($row, 5).Value = $codreg
    $ws.Cells.Item($row, 6).Value = $tipo
    $ws.Cells.Item($row, 7).Value = $productoId
    $ws.Cells.Item($row, 8).Value = $producto
    $ws.Cells.Item($row, 9).Value = $categoriaId
    $ws.Cells.Item($row, 10).Value = $categoria
    $ws.Cells.Item($row, 11).Value = $variedad
    $ws.Cells.Item($row, 12).Value = $r.Calidad
    $ws.Cells.Item($row, 13).Value = $r.Volumen
    $ws.Cells.Item($row, 14).Value = $r.Precio
    $ws.Cells.Item($row, 15).Value = $r.Precio
    $ws.Cells.Item($row, 16).Value = $r.Precio
    $ws.Cells.Item($row, 17).Value = $unidad
    $ws.Cells.Item($row, 18).Value = $origen
    $ws.Cells.Item($row, 19).Value = $r.PrecioKg
    $ws.Cells.Item($row, 20).Value = $kgUnidad
}
